$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mutex-stays-unlocked bug: L2 and L3 should be 1 instead of 0
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 1

# Move the active selection to L2, matching the saved cursor position
$ws.Range("L2").Select()
